$wb = $excel.ActiveWorkbook
$wsReg = $wb.Worksheets.Item("registration")
$wsOne = $wb.Worksheets.Item("Sheet1")

# --- registration sheet: D2 gets a (too long) UCI code typed in for validation testing ---
$wsReg.Range("D2").Value = 91200054323

# --- registration sheet: new row 3, junk/validation test data typed by hand ---
$wsReg.Range("A3").Value = 2
$wsReg.Range("B3").Value = "POOd .djsjs "
$wsReg.Range("C3").Value = 'selrect * from "Top"'
$wsReg.Range("D3").Value = 234
$wsReg.Range("E3").Value = "LTU"
$wsReg.Range("F3").Value = "LTU"
$wsReg.Range("G3").Value = 35697
$wsReg.Range("H3").Value = "female"
$wsReg.Range("I3").Value = "women"

# G3 should look like a typed date (same number format as G2)
$wsReg.Range("G2").Copy() | Out-Null
$wsReg.Range("G3").PasteSpecial(-4122) | Out-Null

# --- registration sheet: new row 4, a competitor record copied over from Sheet1 row 11 ---
$wsReg.Range("A4").Value = 10
$wsReg.Range("B4").Value = "JONES"
$wsReg.Range("C4").Value = "Hayley"
$wsReg.Range("D4").Value = 10009084739
$wsReg.Range("E4").Value = "Team Wales"
$wsReg.Range("F4").Value = "GRB"
$wsReg.Range("G4").Value = 35195
$wsReg.Range("H4").Value = "female"
$wsReg.Range("I4").Value = "women"

# Match the formatting of the source row (Sheet1!A11:I11) cell by cell
foreach ($col in @("A", "C", "F", "H", "I")) {
    $wsReg.Range("B1").Copy() | Out-Null
    $wsReg.Range("$col`4").PasteSpecial(-4122) | Out-Null
}
$wsOne.Range("G11").Copy() | Out-Null
$wsReg.Range("G4").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# --- selection left on Sheet1 at the copied row, then back on registration where the user ended up ---
$wsOne.Range("A11:I11").Select() | Out-Null
$wsReg.Range("B8").Select() | Out-Null
